$wb = $excel.ActiveWorkbook

$wsALC = $wb.Worksheets.Item("ALC")
$wsARM = $wb.Worksheets.Item("ARM")
$wsBSM = $wb.Worksheets.Item("BSM")
$wsCRP = $wb.Worksheets.Item("CRP")
$wsCUL = $wb.Worksheets.Item("CUL")
$wsGSM = $wb.Worksheets.Item("GSM")
$wsLTW = $wb.Worksheets.Item("LTW")
$wsWVR = $wb.Worksheets.Item("WVR")

$wsALC.Range("H7").Value = 3000
$wsALC.Range("J7").Value = 5000
$wsALC.Range("L7").Value = 5000
$wsALC.Range("N7").Value = -5224
$wsALC.Range("H13").Value = 800
$wsALC.Range("J13").Value = 797.5
$wsALC.Range("L13").Value = 797.5
$wsALC.Range("N13").Value = -1135.5
$wsALC.Range("H14").Value = 3000
$wsALC.Range("J14").Value = 5000
$wsALC.Range("L14").Value = 5000
$wsALC.Range("N14").Value = -5382
$wsALC.Range("H16").Value = 2012.75
$wsALC.Range("I16").Value = 250.5
$wsALC.Range("J16").Value = 3775
$wsALC.Range("K16").Value = 250.5
$wsALC.Range("L16").Value = 3775
$wsALC.Range("M16").Value = -20.5
$wsALC.Range("N16").Value = -4235
$wsALC.Range("H19").Value = 1375.5
$wsALC.Range("J19").Value = 2502
$wsALC.Range("L19").Value = 2502
$wsALC.Range("N19").Value = -2852
$wsALC.Range("H38").Value = 1073.0714
$wsALC.Range("I38").Value = 585.4167
$wsALC.Range("K38").Value = 1756.2501
$wsALC.Range("M38").Value = -1384.2501
$wsALC.Range("H69").Value = 31998
$wsALC.Range("J69").Value = 26663.334
$wsALC.Range("L69").Value = 79990.00199999999
$wsALC.Range("N69").Value = -81738.00199999999
$wsALC.Range("H70").Value = 3008.5454
$wsALC.Range("I70").Value = 2349.5
$wsALC.Range("J70").Value = 3799.4
$wsALC.Range("K70").Value = 7048.5
$wsALC.Range("L70").Value = 11398.2
$wsALC.Range("M70").Value = -6778.5
$wsALC.Range("N70").Value = -11938.2
$wsALC.Range("H72").Value = 31998
$wsALC.Range("J72").Value = 26663.334
$wsALC.Range("L72").Value = 239970.006
$wsALC.Range("N72").Value = -248706.006
$wsALC.Range("H73").Value = 3008.5454
$wsALC.Range("I73").Value = 2349.5
$wsALC.Range("J73").Value = 3799.4
$wsALC.Range("K73").Value = 7048.5
$wsALC.Range("L73").Value = 11398.2
$wsALC.Range("M73").Value = -6112.5
$wsALC.Range("N73").Value = -13270.2
$wsALC.Range("H74").Value = 7777.778
$wsALC.Range("I74").Value = 7142.857
$wsALC.Range("J74").Value = 10000
$wsALC.Range("K74").Value = 7142.857
$wsALC.Range("L74").Value = 10000
$wsALC.Range("M74").Value = -6206.857
$wsALC.Range("N74").Value = -11872
$wsALC.Range("H77").Value = 7777.778
$wsALC.Range("I77").Value = 7142.857
$wsALC.Range("J77").Value = 10000
$wsALC.Range("K77").Value = 35714.285
$wsALC.Range("L77").Value = 50000
$wsALC.Range("M77").Value = -31034.285
$wsALC.Range("N77").Value = -59360
$wsALC.Range("H100").Value = 1955.2778
$wsALC.Range("I100").Value = 1684.4615
$wsALC.Range("K100").Value = 1684.4615
$wsALC.Range("M100").Value = -1143.4615
$wsALC.Range("H111").Value = 936.3333
$wsALC.Range("I111").Value = 936.3333
$wsALC.Range("K111").Value = 2808.9999
$wsALC.Range("M111").Value = 258.0001000000002
$wsALC.Range("H125").Value = 2032
$wsALC.Range("I125").Value = 2032
$wsALC.Range("K125").Value = 18288
$wsALC.Range("M125").Value = -15828
$wsALC.Range("H137").Value = 8811.857
$wsALC.Range("J137").Value = 14918.637
$wsALC.Range("L137").Value = 44755.911
$wsALC.Range("N137").Value = -49855.911
$wsALC.Range("H138").Value = 2269.806
$wsALC.Range("I138").Value = 1378.0834
$wsALC.Range("K138").Value = 4134.2502
$wsALC.Range("M138").Value = 1005.7498
$wsARM.Range("H3").Value = 2450
$wsARM.Range("I3").Value = 2400
$wsARM.Range("J3").Value = 2500
$wsARM.Range("K3").Value = 2400
$wsARM.Range("L3").Value = 2500
$wsARM.Range("M3").Value = -2285
$wsARM.Range("N3").Value = -2730
$wsARM.Range("H4").Value = 200
$wsARM.Range("I4").Value = 200
$wsARM.Range("K4").Value = 200
$wsARM.Range("M4").Value = -84
$wsARM.Range("H32").Value = 155583.11
$wsARM.Range("I32").Value = 160062.1
$wsARM.Range("J32").Value = 12255.5
$wsARM.Range("K32").Value = 160062.1
$wsARM.Range("L32").Value = 12255.5
$wsARM.Range("M32").Value = -159775.1
$wsARM.Range("N32").Value = -12829.5
$wsARM.Range("H34").Value = 300000
$wsARM.Range("I34").Value = 300000
$wsARM.Range("K34").Value = 300000
$wsARM.Range("M34").Value = -299729
$wsARM.Range("H38").Value = 8673.333000000001
$wsARM.Range("I38").Value = 0
$wsARM.Range("J38").Value = 8673.333000000001
$wsARM.Range("K38").Value = 0
$wsARM.Range("L38").Value = 8673.333000000001
$wsARM.Range("M38").Value = $null
$wsARM.Range("N38").Value = -9607.333000000001
$wsARM.Range("H43").Value = 162252.22
$wsARM.Range("J43").Value = 239474.2
$wsARM.Range("L43").Value = 239474.2
$wsARM.Range("N43").Value = -240100.2
$wsARM.Range("H45").Value = 4424.75
$wsARM.Range("I45").Value = 3899.6667
$wsARM.Range("K45").Value = 3899.6667
$wsARM.Range("M45").Value = -3522.6667
$wsARM.Range("H61").Value = 6805.5
$wsARM.Range("I61").Value = 6662.4546
$wsARM.Range("K61").Value = 6662.4546
$wsARM.Range("M61").Value = -6450.4546
$wsARM.Range("H74").Value = 5814.614
$wsARM.Range("I74").Value = 909.3333
$wsARM.Range("K74").Value = 909.3333
$wsARM.Range("M74").Value = -35.33330000000001
$wsARM.Range("H77").Value = 5814.614
$wsARM.Range("I77").Value = 909.3333
$wsARM.Range("K77").Value = 4546.6665
$wsARM.Range("M77").Value = -178.6665000000003
$wsARM.Range("H102").Value = 10000
$wsARM.Range("I102").Value = 10000
$wsARM.Range("K102").Value = 10000
$wsARM.Range("M102").Value = -8378
$wsARM.Range("H122").Value = 2037.1
$wsARM.Range("I122").Value = 2005.2222
$wsARM.Range("J122").Value = 2324
$wsARM.Range("K122").Value = 6015.6666
$wsARM.Range("L122").Value = 6972
$wsARM.Range("M122").Value = -3565.6666
$wsARM.Range("N122").Value = -11872
$wsARM.Range("H132").Value = 1193244.2
$wsARM.Range("I132").Value = 1318515.6
$wsARM.Range("J132").Value = 3166.5
$wsARM.Range("K132").Value = 3955546.8
$wsARM.Range("L132").Value = 9499.5
$wsARM.Range("M132").Value = -3953016.8
$wsARM.Range("N132").Value = -14559.5
$wsARM.Range("H133").Value = 62257.855
$wsARM.Range("J133").Value = 62257.855
$wsARM.Range("L133").Value = 62257.855
$wsARM.Range("N133").Value = -67317.85500000001
$wsARM.Range("H136").Value = 6805.5
$wsARM.Range("I136").Value = 6662.4546
$wsARM.Range("K136").Value = 19987.3638
$wsARM.Range("M136").Value = -17437.3638
$wsBSM.Range("H7").Value = 1167500
$wsBSM.Range("J7").Value = 2400
$wsBSM.Range("L7").Value = 2400
$wsBSM.Range("N7").Value = -2626
$wsBSM.Range("I10").Value = 0
$wsBSM.Range("J10").Value = 4250
$wsBSM.Range("K10").Value = 0
$wsBSM.Range("L10").Value = 4250
$wsBSM.Range("M10").Value = $null
$wsBSM.Range("N10").Value = -4530
$wsBSM.Range("H20").Value = 79683.08
$wsBSM.Range("I20").Value = 146819
$wsBSM.Range("K20").Value = 146819
$wsBSM.Range("M20").Value = -146572
$wsBSM.Range("H22").Value = 450
$wsBSM.Range("I22").Value = 300
$wsBSM.Range("K22").Value = 300
$wsBSM.Range("M22").Value = -127
$wsBSM.Range("H86").Value = 1152.9048
$wsBSM.Range("I86").Value = 1200.5294
$wsBSM.Range("J86").Value = 950.5
$wsBSM.Range("K86").Value = 1200.5294
$wsBSM.Range("L86").Value = 950.5
$wsBSM.Range("M86").Value = -77.5293999999999
$wsBSM.Range("N86").Value = -3196.5
$wsBSM.Range("H89").Value = 1152.9048
$wsBSM.Range("I89").Value = 1200.5294
$wsBSM.Range("J89").Value = 950.5
$wsBSM.Range("K89").Value = 6002.646999999999
$wsBSM.Range("L89").Value = 4752.5
$wsBSM.Range("M89").Value = -386.646999999999
$wsBSM.Range("N89").Value = -15984.5
$wsBSM.Range("H94").Value = 3122.6
$wsBSM.Range("I94").Value = 3673.3333
$wsBSM.Range("J94").Value = 919.6667
$wsBSM.Range("K94").Value = 3673.3333
$wsBSM.Range("L94").Value = 919.6667
$wsBSM.Range("M94").Value = -3222.3333
$wsBSM.Range("N94").Value = -1821.6667
$wsBSM.Range("H107").Value = 0
$wsBSM.Range("I107").Value = 0
$wsBSM.Range("J107").Value = 0
$wsBSM.Range("K107").Value = 0
$wsBSM.Range("L107").Value = 0
$wsBSM.Range("M107").Value = $null
$wsBSM.Range("N107").Value = $null
$wsBSM.Range("H134").Value = 4779.6343
$wsBSM.Range("I134").Value = 2234.3547
$wsBSM.Range("J134").Value = 12670
$wsBSM.Range("K134").Value = 6703.0641
$wsBSM.Range("L134").Value = 38010
$wsBSM.Range("M134").Value = -4168.0641
$wsBSM.Range("N134").Value = -43080
$wsCRP.Range("H31").Value = 2570.65
$wsCRP.Range("J31").Value = 2778.111
$wsCRP.Range("L31").Value = 2778.111
$wsCRP.Range("N31").Value = -3368.111
$wsCRP.Range("H34").Value = 2570.65
$wsCRP.Range("J34").Value = 2778.111
$wsCRP.Range("L34").Value = 2778.111
$wsCRP.Range("N34").Value = -3182.111
$wsCRP.Range("H35").Value = 2827.1667
$wsCRP.Range("J35").Value = 3398.3333
$wsCRP.Range("L35").Value = 3398.3333
$wsCRP.Range("N35").Value = -3986.3333
$wsCRP.Range("H58").Value = 24149.834
$wsCRP.Range("I58").Value = 8724.75
$wsCRP.Range("J58").Value = 55000
$wsCRP.Range("K58").Value = 8724.75
$wsCRP.Range("L58").Value = 55000
$wsCRP.Range("M58").Value = -8521.75
$wsCRP.Range("N58").Value = -55406
$wsCRP.Range("H68").Value = 51750.5
$wsCRP.Range("J68").Value = 51750.5
$wsCRP.Range("L68").Value = 51750.5
$wsCRP.Range("N68").Value = -53248.5
$wsCRP.Range("H71").Value = 51750.5
$wsCRP.Range("J71").Value = 51750.5
$wsCRP.Range("L71").Value = 155251.5
$wsCRP.Range("N71").Value = -162739.5
$wsCRP.Range("H132").Value = 2402.7778
$wsCRP.Range("I132").Value = 2406
$wsCRP.Range("K132").Value = 7218
$wsCRP.Range("M132").Value = -4688
$wsCRP.Range("H136").Value = 24149.834
$wsCRP.Range("I136").Value = 8724.75
$wsCRP.Range("J136").Value = 55000
$wsCRP.Range("K136").Value = 26174.25
$wsCRP.Range("L136").Value = 165000
$wsCRP.Range("M136").Value = -23624.25
$wsCRP.Range("N136").Value = -170100
$wsCUL.Range("H4").Value = 7370340
$wsCUL.Range("I4").Value = 12727604
$wsCUL.Range("J4").Value = 4102.25
$wsCUL.Range("K4").Value = 38182812
$wsCUL.Range("L4").Value = 12306.75
$wsCUL.Range("M4").Value = -38182700
$wsCUL.Range("N4").Value = -12530.75
$wsCUL.Range("H5").Value = 1293.7805
$wsCUL.Range("I5").Value = 1179.1428
$wsCUL.Range("K5").Value = 3537.4284
$wsCUL.Range("M5").Value = -3425.4284
$wsCUL.Range("H11").Value = 100000140
$wsCUL.Range("I11").Value = 78
$wsCUL.Range("K11").Value = 234
$wsCUL.Range("M11").Value = -94
$wsCUL.Range("H95").Value = 0
$wsCUL.Range("J95").Value = 0
$wsCUL.Range("L95").Value = 0
$wsCUL.Range("N95").Value = $null
$wsCUL.Range("H105").Value = 22166.666
$wsCUL.Range("J105").Value = 22166.666
$wsCUL.Range("L105").Value = 66499.99800000001
$wsCUL.Range("N105").Value = -71741.99800000001
$wsCUL.Range("H116").Value = 5189
$wsCUL.Range("I116").Value = 4387.3335
$wsCUL.Range("K116").Value = 13162.0005
$wsCUL.Range("M116").Value = -9720.000499999998
$wsCUL.Range("H135").Value = 1293.7805
$wsCUL.Range("I135").Value = 1179.1428
$wsCUL.Range("K135").Value = 10612.2852
$wsCUL.Range("M135").Value = -8077.2852
$wsGSM.Range("H17").Value = 2533.3333
$wsGSM.Range("I17").Value = 100
$wsGSM.Range("J17").Value = 3750
$wsGSM.Range("K17").Value = 100
$wsGSM.Range("L17").Value = 3750
$wsGSM.Range("M17").Value = 68
$wsGSM.Range("N17").Value = -4086
$wsGSM.Range("H52").Value = 0
$wsGSM.Range("J52").Value = 0
$wsGSM.Range("L52").Value = 0
$wsGSM.Range("N52").Value = $null
$wsGSM.Range("H102").Value = 1276.2439
$wsGSM.Range("I102").Value = 1276.2439
$wsGSM.Range("K102").Value = 1276.2439
$wsGSM.Range("M102").Value = 345.7561000000001
$wsGSM.Range("H122").Value = 43032.6
$wsGSM.Range("I122").Value = 52815.85
$wsGSM.Range("J122").Value = 3899.6
$wsGSM.Range("K122").Value = 158447.55
$wsGSM.Range("L122").Value = 11698.8
$wsGSM.Range("M122").Value = -155997.55
$wsGSM.Range("N122").Value = -16598.8
$wsGSM.Range("H125").Value = 20197.334
$wsGSM.Range("I125").Value = 20296
$wsGSM.Range("J125").Value = 20000
$wsGSM.Range("K125").Value = 20296
$wsGSM.Range("L125").Value = 20000
$wsGSM.Range("M125").Value = -17836
$wsGSM.Range("N125").Value = -24920
$wsGSM.Range("H132").Value = 10959.125
$wsGSM.Range("I132").Value = 13042.774
$wsGSM.Range("J132").Value = 3782.111
$wsGSM.Range("K132").Value = 39128.322
$wsGSM.Range("L132").Value = 11346.333
$wsGSM.Range("M132").Value = -36598.322
$wsGSM.Range("N132").Value = -16406.333
$wsLTW.Range("H7").Value = 4567.909
$wsLTW.Range("I7").Value = 3608.1428
$wsLTW.Range("K7").Value = 3608.1428
$wsLTW.Range("M7").Value = -3496.1428
$wsLTW.Range("H40").Value = 2617.5217
$wsLTW.Range("I40").Value = 1850.8125
$wsLTW.Range("K40").Value = 1850.8125
$wsLTW.Range("M40").Value = -1714.8125
$wsLTW.Range("H46").Value = 4649.2856
$wsLTW.Range("I46").Value = 1026
$wsLTW.Range("J46").Value = 5030.684
$wsLTW.Range("K46").Value = 1026
$wsLTW.Range("L46").Value = 5030.684
$wsLTW.Range("M46").Value = -838
$wsLTW.Range("N46").Value = -5406.684
$wsLTW.Range("H68").Value = 2000
$wsLTW.Range("I68").Value = 0
$wsLTW.Range("K68").Value = 0
$wsLTW.Range("M68").Value = $null
$wsLTW.Range("H71").Value = 2000
$wsLTW.Range("I71").Value = 0
$wsLTW.Range("K71").Value = 0
$wsLTW.Range("M71").Value = $null
$wsLTW.Range("H80").Value = 0
$wsLTW.Range("J80").Value = 0
$wsLTW.Range("L80").Value = 0
$wsLTW.Range("N80").Value = $null
$wsLTW.Range("H83").Value = 0
$wsLTW.Range("J83").Value = 0
$wsLTW.Range("L83").Value = 0
$wsLTW.Range("N83").Value = $null
$wsLTW.Range("H92").Value = 50000
$wsLTW.Range("J92").Value = 50000
$wsLTW.Range("L92").Value = 50000
$wsLTW.Range("N92").Value = -54992
$wsLTW.Range("H122").Value = 4172.029
$wsLTW.Range("I122").Value = 3260.84
$wsLTW.Range("K122").Value = 9782.52
$wsLTW.Range("M122").Value = -7332.52
$wsLTW.Range("H126").Value = 4567.909
$wsLTW.Range("I126").Value = 3608.1428
$wsLTW.Range("K126").Value = 10824.4284
$wsLTW.Range("M126").Value = -8354.428400000001
$wsLTW.Range("H132").Value = 2780277.8
$wsLTW.Range("I132").Value = 3848388.5
$wsLTW.Range("K132").Value = 11545165.5
$wsLTW.Range("M132").Value = -11542635.5
$wsWVR.Range("H80").Value = 0
$wsWVR.Range("J80").Value = 0
$wsWVR.Range("L80").Value = 0
$wsWVR.Range("N80").Value = $null
$wsWVR.Range("H81").Value = 3263.2222
$wsWVR.Range("I81").Value = 2083.5
$wsWVR.Range("J81").Value = 5622.6665
$wsWVR.Range("K81").Value = 4167
$wsWVR.Range("L81").Value = 11245.333
$wsWVR.Range("M81").Value = -3106
$wsWVR.Range("N81").Value = -13367.333
$wsWVR.Range("H83").Value = 0
$wsWVR.Range("J83").Value = 0
$wsWVR.Range("L83").Value = 0
$wsWVR.Range("N83").Value = $null
$wsWVR.Range("H84").Value = 3263.2222
$wsWVR.Range("I84").Value = 2083.5
$wsWVR.Range("J84").Value = 5622.6665
$wsWVR.Range("K84").Value = 20835
$wsWVR.Range("L84").Value = 56226.665
$wsWVR.Range("M84").Value = -15531
$wsWVR.Range("N84").Value = -66834.66500000001
$wsWVR.Range("H107").Value = 1356.075
$wsWVR.Range("I107").Value = 850.4706
$wsWVR.Range("J107").Value = 4221.1665
$wsWVR.Range("K107").Value = 2551.4118
$wsWVR.Range("L107").Value = 12663.4995
$wsWVR.Range("M107").Value = -631.4117999999999
$wsWVR.Range("N107").Value = -16503.4995
$wsWVR.Range("H122").Value = 38272.03
$wsWVR.Range("I122").Value = 1905.65
$wsWVR.Range("K122").Value = 5716.950000000001
$wsWVR.Range("M122").Value = -3266.950000000001
$wsWVR.Range("H126").Value = 1859.4231
$wsWVR.Range("I126").Value = 1887.8572
$wsWVR.Range("K126").Value = 5663.571599999999
$wsWVR.Range("M126").Value = -3193.571599999999
$wsWVR.Range("H136").Value = 1759.0625
$wsWVR.Range("I136").Value = 1812.9166
$wsWVR.Range("J136").Value = 1597.5
$wsWVR.Range("K136").Value = 5438.7498
$wsWVR.Range("L136").Value = 4792.5
$wsWVR.Range("M136").Value = -2888.7498
$wsWVR.Range("N136").Value = -9892.5
